$wb = $excel.ActiveWorkbook

# --- Rename "Sheet3" to "mbook_event" and populate it as the mbook_events table schema ---
$ws1 = $wb.Worksheets.Item("Sheet3")
$ws1.Name = "mbook_event"

$mbookEvent = @(
    @("mbook_events", $null),
    @("id", "int"),
    @("event_id", "nvarchar(50)"),
    @("sport", "nvarchar(500)"),
    @("country", "nvarchar(500)"),
    @("competition", "nvarchar(500)"),
    @("start_time", "nvarchar(500)"),
    @("home", "nvarchar(500)"),
    @("away", "nvarchar(500)")
)

for ($i = 0; $i -lt $mbookEvent.Length; $i++) {
    $row = $i + 1
    $ws1.Cells.Item($row, 1).Value = $mbookEvent[$i][0]
    if ($null -ne $mbookEvent[$i][1]) {
        $ws1.Cells.Item($row, 2).Value = $mbookEvent[$i][1]
    }
}

$ws1.Columns.Item(1).ColumnWidth = 11.85546875
$ws1.Columns.Item(2).ColumnWidth = 13.140625

# --- Add new sheet "mbook_market" after mbook_event, populate as mbook_market table schema ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "mbook_market"

$mbookMarket = @(
    @("mbook_market", $null),
    @("id", "int"),
    @("timespan", "int"),
    @("event_id", "nvarchar(50)"),
    @("market_name", "nvarchar(500)"),
    @("runner_name", "nvarchar(500)"),
    @("depth_no", "nvarchar(50)"),
    @("type", "nvarchar(50)"),
    @("odd", "nvarchar(50)"),
    @("amount", "nvarchar(50)")
)

for ($i = 0; $i -lt $mbookMarket.Length; $i++) {
    $row = $i + 1
    $ws2.Cells.Item($row, 1).Value = $mbookMarket[$i][0]
    if ($null -ne $mbookMarket[$i][1]) {
        $ws2.Cells.Item($row, 2).Value = $mbookMarket[$i][1]
    }
}

$ws2.Columns.Item(1).ColumnWidth = 15.42578125
$ws2.Columns.Item(2).ColumnWidth = 17.85546875

# --- Selections matching the target state ---
$ws1.Range("K18").Select() | Out-Null
$ws2.Range("J19").Select() | Out-Null

# --- mbook_market is the active/selected tab ---
$ws2.Select()

$wb.Saved = $false
